# "Generate Report for Handback"
#
# The localization-status workbook gets refreshed once the zh-cn and de-de
# handback packages are produced: the Overview status text moves from
# "Ready for handoff" to "Handed back: in sync with en-US", and each
# language sheet's "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" row gets filled in (plus a hyperlink on the
# newly-populated target-file cell). Column widths are widened a bit to
# fit the new, longer values.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"
$mdFile    = "6c82ee76-023d-4634-b535-6d77c23f1aae.md"
$mdUrl     = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/016688299d9b721b215097c82d80840a9bc96607/e2e/6c82ee76-023d-4634-b535-6d77c23f1aae.md"

# ---------------------------------------------------------------------
# Overview sheet: status cells for zh-cn / de-de now read "Handed back…"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# Widen the (now longer) status columns to fit the new text.
$wsOverview.Range("E1").ColumnWidth = 29.17
$wsOverview.Range("F1").ColumnWidth = 29.17

# ---------------------------------------------------------------------
# zh-cn sheet: record the handback target/handback files + datetime
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $mdUrl, "", "", $mdFile)
$wsZhCn.Range("I2").Value = $mdFile
$wsZhCn.Range("J2").Value = "6c82ee76-023d-4634-b535-6d77c23f1aae.1e98465918153e8e6845228c4b6ef77711dd9936.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-28 00:56:27"

$wsZhCn.Range("C1").ColumnWidth = 29.17
$wsZhCn.Range("I1").ColumnWidth = 39.17
$wsZhCn.Range("J1").ColumnWidth = 39.17

# ---------------------------------------------------------------------
# de-de sheet: record the handback target/handback files + datetime
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $mdUrl, "", "", $mdFile)
$wsDeDe.Range("I2").Value = $mdFile
$wsDeDe.Range("J2").Value = "6c82ee76-023d-4634-b535-6d77c23f1aae.1e98465918153e8e6845228c4b6ef77711dd9936.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-28 00:56:33"

$wsDeDe.Range("C1").ColumnWidth = 29.17
$wsDeDe.Range("I1").ColumnWidth = 39.17
$wsDeDe.Range("J1").ColumnWidth = 39.17
